$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.866.49"
$ws.Range("E2").Value = "  -1.08%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.37"
$ws.Range("E3").Value = "  -0.84%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("E5").Value = "  -4.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.72"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("E7").Value = "  -0.24%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3152"
$ws.Range("E8").Value = "  -3.39%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.37"
$ws.Range("E9").Value = "  -5.67%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07217"
$ws.Range("E10").Value = "  +2.25%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08097"
$ws.Range("E11").Value = "  -0.03%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7667"
$ws.Range("E12").Value = "  +0.02%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.928.29"
$ws.Range("E13").Value = "  +0.90%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.500"
$ws.Range("E14").Value = "  +4.07%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.53"
$ws.Range("E15").Value = "  -0.40%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.166"
$ws.Range("E16").Value = "  +4.47%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.887.23"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.97"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.05"
$ws.Range("E19").Value = "  -0.58%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007804"
$ws.Range("E20").Value = "  +0.12%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.154.45"
$ws.Range("E21").Value = "  -0.43%  "

# Row 22
$ws.Range("E22").Value = "  -0.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.129"
$ws.Range("E23").Value = "  +15.42%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.16%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1648"
$ws.Range("E25").Value = "  -1.71%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.433"
$ws.Range("E26").Value = "  +1.25%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.31"
$ws.Range("E27").Value = "  -2.08%  "

# Row 28
$ws.Range("E28").Value = "  -1.40%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.060"
$ws.Range("E29").Value = "  -2.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.407"
$ws.Range("E30").Value = "  +2.45%  "

# Row 31
$ws.Range("E31").Value = "  +1.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.507"
$ws.Range("E32").Value = "  +4.65%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.125"
$ws.Range("E33").Value = "  +0.81%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05576"
$ws.Range("E34").Value = "  -6.37%  "

# Row 35
$ws.Range("E35").Value = "  -0.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7440"
$ws.Range("E36").Value = "  +1.29%  "

# Row 37
$ws.Range("E37").Value = "  -0.40%  "

# Row 38
$ws.Range("E38").Value = "  -3.82%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01925"
$ws.Range("E39").Value = "  -0.25%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.775"
$ws.Range("E40").Value = "  -0.71%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.145.82"
$ws.Range("E41").Value = "  +13.84%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.27"
$ws.Range("E42").Value = "  +1.43%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4430"
$ws.Range("E43").Value = "  -0.84%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.866"
$ws.Range("E44").Value = "  -1.75%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8508"
$ws.Range("E45").Value = "  -0.40%  "

# Row 46
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.04"
$ws.Range("E46").Value = "  +1.60%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.30%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.886"
$ws.Range("E48").Value = "  -1.15%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.995"
$ws.Range("E49").Value = "  +1.56%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.480"
$ws.Range("E50").Value = "  -1.48%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.022"
$ws.Range("E51").Value = "  +11.04%  "
